$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "Laboratório"
$ws.Range("A15").Value = "UATS"

$ws.Range("A13").Copy()
$ws.Range("A14:A15").PasteSpecial(-4122)

$ws.Range("A25").Select()
